# Season record columns: the scraper previously only pulled team/player
# statistics, not the team's season record (Wins/Losses/Ties). This adds
# those three columns to the right of the existing stat table and fills
# in the 2008 Mets' record (89-73-0) for every player row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Figure out where the existing table ends so this keeps working even if
# the sheet layout shifts (last used row/column of the current data).
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

$winsCol   = $lastCol + 1
$lossesCol = $lastCol + 2
$tiesCol   = $lastCol + 3

# New header cells should look like the rest of row 1 (bold, centered,
# thin border) -- copy the formatting from the last existing header cell
# before writing the new header text into place.
$lastHeaderCell = $ws.Cells.Item(1, $lastCol)
$newHeaderRange = $ws.Range($ws.Cells.Item(1, $winsCol), $ws.Cells.Item(1, $tiesCol))
$lastHeaderCell.Copy()
$newHeaderRange.PasteSpecial(-4122)

$ws.Cells.Item(1, $winsCol).Value   = "Wins"
$ws.Cells.Item(1, $lossesCol).Value = "Losses"
$ws.Cells.Item(1, $tiesCol).Value   = "Ties"

# The 2008 New York Mets finished the season 89-73-0.
$wins = 89
$losses = 73
$ties = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $winsCol).Value   = $wins
    $ws.Cells.Item($r, $lossesCol).Value = $losses
    $ws.Cells.Item($r, $tiesCol).Value   = $ties
}
